# chore: replace gui_element Slider by Spinbox, remove gui_element SimpleText
# for intervals, et al. (DEV-2501) (#452)
#
# Row 4 of the sheet described the "hasInterval" property. It is turned into
# a "hasInteger" property (IntValue / Spinbox gui element) instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "hasInteger"
$ws.Range("B4").Value = "has Integer"
$ws.Range("C4").Value = "Zahl"
$ws.Range("G4").Value = "integer"
$ws.Range("H4").Value = "zahl"
$ws.Range("L4").Value = "hasValue"
$ws.Range("M4").Value = "IntValue"
$ws.Range("N4").Value = "Spinbox"
$ws.Range("O4").Value = "max: 10, min: 5, rows: 10"
